$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "dgsdgdsvgdsgfd"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 3
